# Updated with best results
# - Removes the blank spacer row 9 (the "Integer" section previously had a
#   2-row gap before "Floating Point"; every other section only had 1 row).
#   Deleting row 9 shifts everything below it up by one row, matching the
#   other sections' single-row gap.
# - Refreshes a batch of result values across the four sections with newer
#   measurements.
# - Restores the active-cell selection to F14.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the extra blank row that used to separate the "Integer" block from
# the "Floating Point" block; this shifts rows 10-30 up to 9-29.
$ws.Rows.Item(9).Delete()

$updates = [ordered]@{
    "C3"  = 12.778
    "L3"  = 12.39
    "C4"  = 9.9309999999999992
    "I4"  = 7.9870000000000001
    "L4"  = 9.5719999999999992
    "C5"  = 8.4600000000000009
    "L5"  = 7.3739999999999997
    "C6"  = 7.492
    "I6"  = 5.351
    "L6"  = 5.8730000000000002
    "C7"  = 6.7140000000000004
    "F7"  = 7.5529999999999999
    "I7"  = 3.1419999999999999
    "L7"  = 4.8840000000000003
    "C8"  = 5.1950000000000003
    "I8"  = 1.3420000000000001
    "L8"  = 3.1480000000000001
    "C10" = 3.3690000000000002
    "I10" = 3.496
    "L10" = 3.27
    "L11" = 2.8719999999999999
    "I12" = 2.5950000000000002
    "L12" = 2.629
    "C13" = 2.665
    "I13" = 2.282
    "I14" = 2.0510000000000002
    "L14" = 2.4409999999999998
    "C15" = 2.6339999999999999
    "I15" = 1.931
    "L15" = 2.2959999999999998
    "L17" = 10.416
    "I18" = 8.2040000000000006
    "I21" = 7.5229999999999997
    "L21" = 8.8539999999999992
    "C22" = 10.268000000000001
    "F22" = 8.8840000000000003
    "I22" = 7.3630000000000004
    "L22" = 8.4629999999999992
    "C24" = 11.509
    "I24" = 12.6
    "L24" = 11.613
    "C25" = 9.6370000000000005
    "I25" = 10.148999999999999
    "L25" = 9.5289999999999999
    "C26" = 8.3030000000000008
    "I26" = 7.4349999999999996
    "L26" = 7.4169999999999998
    "C27" = 7.0839999999999996
    "L27" = 5.9429999999999996
    "C28" = 6.0419999999999998
    "I28" = 3.3
    "L28" = 4.8410000000000002
    "C29" = 3.9159999999999999
    "I29" = 2.0169999999999999
    "L29" = 2.8130000000000002
}

foreach ($ref in $updates.Keys) {
    $ws.Range($ref).Value = $updates[$ref]
}

# Match the author's final cursor position.
$ws.Range("F14").Select()
